$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing "Week 7" hours for member in row 4 (ckc2)
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 9

# Add the missing "Week 7" column (header + values) to the second summary
# table ("Total Hours Spent Per Task"), rows 35-38, column Q
$ws.Range("Q35").Value = "Week 7"
$ws.Range("Q35").Style = $ws.Range("P35").Style

$ws.Range("P36").Value = 3
$ws.Range("Q36").Value = 4

$ws.Range("Q37").Value = 0

$ws.Range("Q38").Value = 4

# Update selection/view to match the new editing position
$ws.Range("U29").Select()

$wb.RecalculateFull()
